$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G2: reorder "Recorded By" list (System moved to front) ---
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# --- G3: reorder "Recorded By" list ---
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- G4: reorder "Recorded By" list ---
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# --- G5: reorder "Recorded By" list ---
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- G7 / H7: reorder list (add Amera.a.saad) and update student count ---
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("H7").Value = "79/251"

# --- G9: reorder "Recorded By" list ---
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# --- L10: Average Attendance % value update (keep stored as text, like original) ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "24.3%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- G12: reorder "Recorded By" list ---
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"

# --- G15: reorder "Recorded By" list ---
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# --- S15: Average Attendance % value update (keep stored as text, like original) ---
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "24.3%"
$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# --- G28: reorder "Recorded By" list ---
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
